$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 194; this shifts existing rows 194-339 down to 195-340.
$ws.Rows(194).Insert()

# Populate the newly inserted row 194 with the new weekly price record.
$ws.Range("A194").Value = 5
$ws.Range("B194").Value = "Macroferia Regional de Talca"
$ws.Range("C194").Value = "Maule"
$ws.Range("D194").Value = 44741
$ws.Range("E194").Value = 7
$ws.Range("F194").Value = 100114014
$ws.Range("G194").Value = "Betarraga"
$ws.Range("H194").Value = "Sin especificar"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 5000
$ws.Range("K194").Value = 700
$ws.Range("L194").Value = 700
$ws.Range("M194").Value = 700
$ws.Range("N194").Value = "$/paquete 5 unidades"
$ws.Range("O194").Value = "Región del Maule"
$ws.Range("P194").Value = 140
$ws.Range("Q194").Value = 5
$ws.Range("R194").Value = "Hortaliza"
